$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = 10
$ws.Range("G5").Value = 2.2
$ws.Range("I5").Value = 3.8
$ws.Range("J5").Value = 3
$ws.Range("W5").Value = 1.53
$ws.Range("X5").Value = 2.38
$ws.Range("AB5").Value = 9
$ws.Range("AC5").Value = 9.5
$ws.Range("AF5").Value = 34
$ws.Range("AH5").Value = 6
$ws.Range("AM5").Value = 17
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3
$ws.Range("G14").Value = 4.2
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 1.67
$ws.Range("K14").Value = 2.3
$ws.Range("Q14").Value = 1.75
$ws.Range("R14").Value = 2.05
$ws.Range("AG14").Value = 13
$ws.Range("AH14").Value = 8
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("Q15").Value = 1.95
$ws.Range("R15").Value = 1.9
$ws.Range("U15").Value = 3.4
$ws.Range("V15").Value = 1.33
$ws.Range("H16").Value = 3.6
$ws.Range("AJ16").Value = 81
$ws.Range("AR16").Value = 1.85
$ws.Range("AS16").Value = 2
$ws.Range("G17").Value = 2.45
$ws.Range("I17").Value = 2.4
$ws.Range("J17").Value = 3.1
$ws.Range("L17").Value = 3
$ws.Range("S17").Value = 1.93
$ws.Range("T17").Value = 1.88
$ws.Range("G18").Value = 2.35
$ws.Range("H18").Value = 3.25
$ws.Range("I18").Value = 2.82
$ws.Range("J18").Value = 2.87
$ws.Range("K18").Value = 2.12
$ws.Range("L18").Value = 3.35
$ws.Range("O18").Value = 1.2
$ws.Range("P18").Value = 3.6
$ws.Range("Q18").Value = 1.6
$ws.Range("R18").Value = 2.05
$ws.Range("U18").Value = 2.4
$ws.Range("V18").Value = 1.44
$ws.Range("Y18").Value = 1.47
$ws.Range("Z18").Value = 2.32
$ws.Range("AA18").Value = 10.75
$ws.Range("AB18").Value = 14
$ws.Range("AD18").Value = 27
$ws.Range("AE18").Value = 17
$ws.Range("AF18").Value = 21
$ws.Range("AH18").Value = 6.6
$ws.Range("AI18").Value = 11
$ws.Range("AJ18").Value = 37
$ws.Range("AK18").Value = 200
$ws.Range("AM18").Value = 17
$ws.Range("AP18").Value = 22
$ws.Range("G19").Value = 3.75
$ws.Range("I19").Value = 1.85
$ws.Range("K19").Value = 2.18
$ws.Range("L19").Value = 2.42
$ws.Range("Z19").Value = 2.02
$ws.Range("AA19").Value = 12.5
$ws.Range("AB19").Value = 22
$ws.Range("AF19").Value = 32
$ws.Range("AG19").Value = 11.75
$ws.Range("AH19").Value = 7
$ws.Range("AL19").Value = 8
$ws.Range("AM19").Value = 9.5
$ws.Range("AO19").Value = 16
$ws.Range("G20").Value = 1.7
$ws.Range("I20").Value = 4.5
$ws.Range("K20").Value = 2.4
$ws.Range("O20").Value = 1.17
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = 1.57
$ws.Range("R20").Value = 2.35
$ws.Range("S20").Value = 1.98
$ws.Range("T20").Value = 1.88
$ws.Range("U20").Value = 2.38
$ws.Range("V20").Value = 1.53
$ws.Range("AA20").Value = 9.5
$ws.Range("AD20").Value = 15
$ws.Range("AK20").Value = 126
$ws.Range("O22").Value = 1.17
$ws.Range("P22").Value = 5
$ws.Range("Q22").Value = 1.57
$ws.Range("R22").Value = 2.35
$ws.Range("S22").Value = 1.98
$ws.Range("T22").Value = 1.88
$ws.Range("U22").Value = 2.38
$ws.Range("V22").Value = 1.53
$ws.Range("G24").Value = 2.7
$ws.Range("H24").Value = 2.9
$ws.Range("I24").Value = 2.88
$ws.Range("K24").Value = 1.83
$ws.Range("L24").Value = 3.75
$ws.Range("M24").Value = 1.13
$ws.Range("N24").Value = 6
$ws.Range("O24").Value = 1.57
$ws.Range("P24").Value = 2.25
$ws.Range("Q24").Value = 2.88
$ws.Range("R24").Value = 1.4
$ws.Range("U24").Value = 6
$ws.Range("V24").Value = 1.13
$ws.Range("W24").Value = 1.62
$ws.Range("X24").Value = 2.2
$ws.Range("Y24").Value = 2.2
$ws.Range("Z24").Value = 1.62
$ws.Range("AB24").Value = 11
$ws.Range("AE24").Value = 29
$ws.Range("AG24").Value = 6
$ws.Range("AI24").Value = 21
$ws.Range("AJ24").Value = 81
$ws.Range("AN24").Value = 12
$ws.Range("AP24").Value = 29
$ws.Range("AR24").Value = 2.03
$ws.Range("AS24").Value = 1.83
$ws.Range("G25").Value = 3.75
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 4.33
$ws.Range("L25").Value = 2.75
$ws.Range("Q25").Value = 2.08
$ws.Range("R25").Value = 1.73
$ws.Range("Y25").Value = 1.91
$ws.Range("Z25").Value = 1.91
$ws.Range("AG25").Value = 9
$ws.Range("AI25").Value = 15
$ws.Range("AK25").Value = 301
$ws.Range("AL25").Value = 7
$ws.Range("AM25").Value = 9
$ws.Range("G26").Value = 8
$ws.Range("I26").Value = 1.36
$ws.Range("L26").Value = 1.83
$ws.Range("M26").Value = 1.04
$ws.Range("N26").Value = 13
$ws.Range("Y26").Value = 2.05
$ws.Range("Z26").Value = 1.7
$ws.Range("AA26").Value = 19
$ws.Range("AC26").Value = 23
$ws.Range("AD26").Value = 101
$ws.Range("AG26").Value = 12
$ws.Range("AO26").Value = 8.5
$ws.Range("I27").Value = 3.85
$ws.Range("J27").Value = 2.45
$ws.Range("K27").Value = 2.18
$ws.Range("L27").Value = 4.15
$ws.Range("Q27").Value = 1.78
$ws.Range("R27").Value = 1.93
$ws.Range("X27").Value = 2.82
$ws.Range("AA27").Value = 7.7
$ws.Range("AB27").Value = 9
$ws.Range("AD27").Value = 15.5
$ws.Range("AF27").Value = 24
$ws.Range("AL27").Value = 12.5
$ws.Range("AN27").Value = 12.5
$ws.Range("AO27").Value = 55
$ws.Range("AP27").Value = 32
$ws.Range("AQ27").Value = 35
$ws.Range("M28").Value = 1.05
$ws.Range("N28").Value = 11
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 4.2
$ws.Range("I29").Value = 1.45
$ws.Range("J29").Value = 5.5
$ws.Range("L29").Value = 1.91
$ws.Range("O29").Value = 1.11
$ws.Range("P29").Value = 6.5
$ws.Range("Q29").Value = 1.4
$ws.Range("R29").Value = 2.88
$ws.Range("U29").Value = 1.91
$ws.Range("V29").Value = 1.8
$ws.Range("AA29").Value = 23
$ws.Range("AC29").Value = 19
$ws.Range("AD29").Value = 67
$ws.Range("AF29").Value = 34
$ws.Range("AH29").Value = 9
$ws.Range("AI29").Value = 13
$ws.Range("AJ29").Value = 34
$ws.Range("AL29").Value = 12
$ws.Range("AM29").Value = 10
$ws.Range("AN29").Value = 9
$ws.Range("AO29").Value = 12
